# Natmi following Dr Hou advice:
# Recompute the S100a8 -> Tlr4 ligand-receptor table to include the
# "sCs" cluster alongside "ECs" and "FAPs" (sending/receiving), expanding
# the 2-row result set to the full 6-row cross table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "S100a8"
$ws.Range("C2").Value2 = "Tlr4"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 1588.882602
$ws.Range("H2").Value2 = 4766.647806
$ws.Range("I2").Value2 = 0.999783480649261
$ws.Range("J2").Value2 = 0.9997834806492609
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 25.94532566666667
$ws.Range("N2").Value2 = 77.835977
$ws.Range("O2").Value2 = 0.5401813355606462
$ws.Range("P2").Value2 = 0.5401813355606462
$ws.Range("Q2").Value2 = 41224.07655499072
$ws.Range("R2").Value2 = 371016.6889949164
$ws.Range("S2").Value2 = 0.5400643758485892
$ws.Range("T2").Value2 = 0.5400643758485892

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "S100a8"
$ws.Range("C3").Value2 = "Tlr4"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 1588.882602
$ws.Range("H3").Value2 = 4766.647806
$ws.Range("I3").Value2 = 0.999783480649261
$ws.Range("J3").Value2 = 0.9997834806492609
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 18.51427066666667
$ws.Range("N3").Value2 = 55.542812
$ws.Range("O3").Value2 = 0.3854668692210787
$ws.Range("P3").Value2 = 0.3854668692210786
$ws.Range("Q3").Value2 = 29417.00255098561
$ws.Range("R3").Value2 = 264753.0229588705
$ws.Range("S3").Value2 = 0.3853834081848236
$ws.Range("T3").Value2 = 0.3853834081848235

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "S100a8"
$ws.Range("C4").Value2 = "Tlr4"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1588.882602
$ws.Range("H4").Value2 = 4766.647806
$ws.Range("I4").Value2 = 0.999783480649261
$ws.Range("J4").Value2 = 0.9997834806492609
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 3.571174000000001
$ws.Range("N4").Value2 = 10.713522
$ws.Range("O4").Value2 = 0.07435179521827505
$ws.Range("P4").Value2 = 0.07435179521827504
$ws.Range("Q4").Value2 = 5674.176237314749
$ws.Range("R4").Value2 = 51067.58613583274
$ws.Range("S4").Value2 = 0.07433569661584812
$ws.Range("T4").Value2 = 0.07433569661584809

# Row 5
$ws.Range("A5").Value2 = "sCs"
$ws.Range("B5").Value2 = "S100a8"
$ws.Range("C5").Value2 = "Tlr4"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.3440983333333333
$ws.Range("H5").Value2 = 1.032295
$ws.Range("I5").Value2 = 0.0002165193507390483
$ws.Range("J5").Value2 = 0.0002165193507390482
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 25.94532566666667
$ws.Range("N5").Value2 = 77.835977
$ws.Range("O5").Value2 = 0.5401813355606462
$ws.Range("P5").Value2 = 0.5401813355606462
$ws.Range("Q5").Value2 = 8.927743319690554
$ws.Range("R5").Value2 = 80.349689877215
$ws.Range("S5").Value2 = 0.0001169597120569431
$ws.Range("T5").Value2 = 0.0001169597120569431

# Row 6
$ws.Range("A6").Value2 = "sCs"
$ws.Range("B6").Value2 = "S100a8"
$ws.Range("C6").Value2 = "Tlr4"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.3440983333333333
$ws.Range("H6").Value2 = 1.032295
$ws.Range("I6").Value2 = 0.0002165193507390483
$ws.Range("J6").Value2 = 0.0002165193507390482
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 18.51427066666667
$ws.Range("N6").Value2 = 55.542812
$ws.Range("O6").Value2 = 0.3854668692210787
$ws.Range("P6").Value2 = 0.3854668692210786
$ws.Range("Q6").Value2 = 6.370729679282223
$ws.Range("R6").Value2 = 57.33656711354001
$ws.Range("S6").Value2 = 0.00008346103625516159
$ws.Range("T6").Value2 = 0.00008346103625516157

# Row 7
$ws.Range("A7").Value2 = "sCs"
$ws.Range("B7").Value2 = "S100a8"
$ws.Range("C7").Value2 = "Tlr4"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.3440983333333333
$ws.Range("H7").Value2 = 1.032295
$ws.Range("I7").Value2 = 0.0002165193507390483
$ws.Range("J7").Value2 = 0.0002165193507390482
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 3.571174000000001
$ws.Range("N7").Value2 = 10.713522
$ws.Range("O7").Value2 = 0.07435179521827505
$ws.Range("P7").Value2 = 0.07435179521827504
$ws.Range("Q7").Value2 = 1.228835021443333
$ws.Range("R7").Value2 = 11.05951519299
$ws.Range("S7").Value2 = 0.00001609860242694359
$ws.Range("T7").Value2 = 0.00001609860242694358
